$wb = $excel.ActiveWorkbook

# --- OFF sheet: update Row 3 (R) values ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 215
$wsOff.Range("C3").Value = 134
$wsOff.Range("D3").Value = 49
$wsOff.Range("E3").Value = 17
$wsOff.Range("F3").Value = 6

# --- DEF sheet: update Row 3 (R) values ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 266
$wsDef.Range("C3").Value = 183
$wsDef.Range("D3").Value = 59
$wsDef.Range("E3").Value = 26
